# Fixing errors in example upload files.
$wb = $excel.ActiveWorkbook

# --- Practitioners sheet: add a missing data row (row 6) ---
$practitioners = $wb.Worksheets.Item("Practitioners")

$practitioners.Cells.Item(6, 1).Value = "PHN999:NFP02"
$practitioners.Cells.Item(6, 2).Value = "P01"
$practitioners.Cells.Item(6, 3).Value = 8
$practitioners.Cells.Item(6, 4).Value = 1
$practitioners.Cells.Item(6, 5).Value = 1973
$practitioners.Cells.Item(6, 6).Value = 2
$practitioners.Cells.Item(6, 7).Value = 1
$practitioners.Cells.Item(6, 8).Value = 1
$practitioners.Cells.Item(6, 9).Value = "tag1"

# Column widths on Practitioners sheet (values chosen so the engine's
# pixel-rounded stored width lands on the target: 14.6640625, 13, 12.83203125)
$practitioners.Columns.Item(1).ColumnWidth = 13.8335
$practitioners.Columns.Item(3).ColumnWidth = 12.167
$practitioners.Columns.Item(6).ColumnWidth = 12.0005

# Selection / active cell on Practitioners sheet
$practitioners.Range("G1:G1048576").Select() | Out-Null

# --- Service Contacts sheet: column width + selection fix ---
$serviceContacts = $wb.Worksheets.Item("Service Contacts")
$serviceContacts.Columns.Item(1).ColumnWidth = 13.667
$serviceContacts.Range("D3").Select() | Out-Null
